$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data order/values for rows 2-14 (header in row 1 stays unchanged)
$data = @(
    @("TIMOTEO BAYONA SHARYN LISSETH", 114),
    @("ZAPATA ZETA ROSA ARACELI", 112),
    @("PEREZ VEGA ANA YSABEL", 110),
    @("GARAVITO LEON IVONNE LISSETH", 99),
    @("PANTA MONZON SHIRLEY MARIBEL", 89),
    @("NIÑO GUERRERO ANYELA MELINA", 80),
    @("CASTRO JUAREZ MARIA ISABEL", 79),
    @("MORENO PALACIOS DAMARIS VANESA", 74),
    @("VALLE SILVA SUTMMER ORFELINDA", 71),
    @("TIZON NUÑEZ FRESIA YAMILI", 70),
    @("CHERO JUAREZ ANYELA TATIANA", 61),
    @("71050834", 1),
    @("MORETO ESPINOZA JUAN ALBERTO", 1)
)

$r = 2
foreach ($row in $data) {
    $cellA = $ws.Cells.Item($r, 1)
    if ($row[0] -eq "71050834") {
        $cellA.NumberFormat = "@"
    }
    $cellA.Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
